$d = $word.ActiveDocument

# 1) Expand the underscores in the "Copies served by Dep. Clerk" line.
$d.Content.Find.Execute(
    "Copies served by Dep. Clerk ___________ on the following date ___________ to:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Copies served by Dep. Clerk ___________________________ on the following date ____________________ to:",
    2
) | Out-Null

# 2) Tighten the first "PS" spacing and add a trailing semicolon to the
#    Prosecutor's Office / Defendant's Attorney / Defendant line.
$d.Content.Find.Execute(
    "Prosecutor’s Office: PS     OM     EM; Defendant’s Attorney: PS     OM     EM; {{ defendant.first_name }} {{ defendant.last_name}}: PS     OM     EM",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Prosecutor’s Office: PS    OM     EM; Defendant’s Attorney: PS     OM     EM; {{ defendant.first_name }} {{ defendant.last_name}}: PS     OM     EM;",
    2
) | Out-Null

# 3) Insert a brand-new paragraph right after that line for the
#    Community Control / County Jail service-of-copies text. Locate the
#    paragraph by its (now updated) content rather than a hard-coded
#    index so the edit is resilient to any other paragraph-count drift.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Prosecutor*Office*PS*OM*EM*") {
        $targetIndex = $i
    }
}

$sourcePara = $d.Paragraphs.Item($targetIndex)
$sourcePara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.InsertBefore("{% if community_control.ordered is true or bond_conditions.monitoring is true %}Community Control: PS    EM;{% endif %}{% if jail_terms.ordered is true or apply_jtc == ‘Sentence’ %}County Jail: PS   EM;{% endif %}")
